$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '303.51'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-4.68%'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '7'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '35.16'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-2.49%'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '7'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.057'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-2.76%'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '7'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07988'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-3.06%'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '7'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.932'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-10.24%'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '7'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '4.056'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-1.97%'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '7'
$ws.Range("B8").Value = 'KuCoinToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '7.747'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-3.22%'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '7'
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.942'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '5.03%'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '7'
$ws.Range("B10").Value = 'MXToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9215'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-0.61%'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '7'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1228'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '21.00%'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '7'
$ws.Range("B12").Value = 'WazirX'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1842'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-2.52%'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '7'
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09362'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '1.90%'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '7'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03567'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-1.58%'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '7'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09850'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.71%'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '7'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001388'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-3.80%'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '7'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.005738'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.27%'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '7'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.499'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '1.25%'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '7'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3447'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '2.13%'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '7'
$ws.Range("B20").Value = 'ProBitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1308'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.56%'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '7'
$ws.Range("B21").Value = 'MCDex'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.034'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.61%'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '7'
$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2464'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '12.48%'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '7'
$ws.Range("B23").Value = 'CoinExToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04498'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-2.18%'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '7'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-2.44%'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '7'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004854'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '2.59%'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '7'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.16%'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '7'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-6.92%'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '7'
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '7'
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '7'
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '7'
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '7'
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '7'
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '7'
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '7'
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '7'
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '7'
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '7'
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '7'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01934'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-3.34%'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '7'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-4.31%'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '7'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007551'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-2.26%'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '7'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009550'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '22.13%'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '7'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1332'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-4.91%'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '7'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002109'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '0.51%'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '7'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01113'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-6.31%'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '7'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006284'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-2.41%'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '7'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '7'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '55.93%'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '7'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-31.42%'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '7'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '7'
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '7'
